$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old Turkish color column (column S) contents; the data is being
# moved (and translated) into column J instead. Use Clear() rather than
# ClearContents() so the cells are fully vacated (not left behind as blank
# used cells), matching the shrunk worksheet dimension in the target file.
$ws.Range("S2").Clear()
$ws.Range("S3").Clear()
$ws.Range("S5").Clear()
$ws.Range("S7").Clear()
$ws.Range("S9").Clear()
$ws.Range("S11").Clear()
$ws.Range("S14").Clear()
$ws.Range("S30").Clear()

# Add the new "wire" column header and translated (English) wire colors in
# column J.
$ws.Range("J2").Value = "blue"
$ws.Range("J3").Value = "blue/white"
$ws.Range("J5").Value = "brown/white"
$ws.Range("J7").Value = "brown"
$ws.Range("J9").Value = "orange/white"
$ws.Range("J11").Value = "orange"
$ws.Range("J14").Value = "green"
$ws.Range("J30").Value = "green/white"
$ws.Range("J1").Value = "wire"

# The wire colors are now the longest strings in column J, so its width
# needs to be recalculated to fit the new content.
$ws.Columns("J:J").AutoFit()

# Update selection to match the final cursor position recorded in the file.
$ws.Range("N35").Select()
